$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 421, pushing rows 421-453
# down to 423-455 (mirrors the weekly refresh: two new price entries are
# prepended and the oldest history keeps its place further down).
$ws.Range("A421:T422").EntireRow.Insert()

# New row 421
$ws.Range("A421").Value = 1
$ws.Range("B421").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C421").Value = "Arica y Parinacota"
$ws.Range("D421").Value = 45265
$ws.Range("E421").Value = 15
$ws.Range("F421").Value = "Fruta"
$ws.Range("G421").Value = 100108
$ws.Range("H421").Value = "Tropicales y subtropicales"
$ws.Range("I421").Value = 100108006
$ws.Range("J421").Value = "Plátano"
$ws.Range("K421").Value = "Sin especificar"
$ws.Range("L421").Value = "Maduro"
$ws.Range("M421").Value = 108
$ws.Range("N421").Value = 18000
$ws.Range("O421").Value = 19000
$ws.Range("P421").Value = 18500
$ws.Range("Q421").Value = "$/caja 20 kilos"
$ws.Range("R421").Value = "Ecuador"
$ws.Range("S421").Value = 925
$ws.Range("T421").Value = 20

# New row 422
$ws.Range("A422").Value = 1
$ws.Range("B422").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C422").Value = "Arica y Parinacota"
$ws.Range("D422").Value = 45265
$ws.Range("E422").Value = 15
$ws.Range("F422").Value = "Fruta"
$ws.Range("G422").Value = 100108
$ws.Range("H422").Value = "Tropicales y subtropicales"
$ws.Range("I422").Value = 100108006
$ws.Range("J422").Value = "Plátano"
$ws.Range("K422").Value = "Sin especificar"
$ws.Range("L422").Value = "Pintón"
$ws.Range("M422").Value = 120
$ws.Range("N422").Value = 20000
$ws.Range("O422").Value = 21000
$ws.Range("P422").Value = 20500
$ws.Range("Q422").Value = "$/caja 20 kilos"
$ws.Range("R422").Value = "Ecuador"
$ws.Range("S422").Value = 1025
$ws.Range("T422").Value = 20
